$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.797.55'
$ws.Range('E2').Value = '  +1.62%  '

$ws.Range('D3').Value = '3.768.85'
$ws.Range('E3').Value = '  -1.07%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '''603.11'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('D6').Value = '''169.47'
$ws.Range('E6').Value = '  -0.52%  '

$ws.Range('D7').Value = '3.767.12'
$ws.Range('E7').Value = '  -1.05%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').Value = '''0.537'
$ws.Range('E9').Value = '  +1.15%  '

$ws.Range('D10').Value = '''0.165'
$ws.Range('E10').Value = '  +3.30%  '

$ws.Range('D11').Value = '''6.36'
$ws.Range('E11').Value = '  +2.76%  '

$ws.Range('D12').Value = '''0.464'
$ws.Range('E12').Value = '  -0.85%  '

$ws.Range('D13').Value = '''38.47'
$ws.Range('E13').Value = '  -0.94%  '

$ws.Range('E14').Value = '  +0.83%  '

$ws.Range('D15').Value = '4.395.92'
$ws.Range('E15').Value = '  -1.05%  '

$ws.Range('D16').Value = '3.753.82'
$ws.Range('E16').Value = '  -1.34%  '

$ws.Range('D17').Value = '68.809.25'
$ws.Range('E17').Value = '  +1.53%  '

$ws.Range('D18').Value = '''7.32'
$ws.Range('E18').Value = '  +0.73%  '

$ws.Range('E19').Value = '  +0.01%  '

$ws.Range('D20').Value = '''17.24'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''496.58'
$ws.Range('E21').Value = '  +0.16%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''10.73'
$ws.Range('E22').Value = '  +11.04%  '

$ws.Range('D23').Value = '''0.732'
$ws.Range('E23').Value = '  -1.82%  '

$ws.Range('D24').Value = '''85.71'
$ws.Range('E24').Value = '  -0.35%  '

$ws.Range('D25').Value = '''0.0000147'
$ws.Range('E25').Value = '  +1.29%  '

$ws.Range('E26').Value = '  -3.18%  '

$ws.Range('D27').Value = '''12.44'
$ws.Range('E27').Value = '  +0.25%  '

$ws.Range('D28').Value = '''10.21'
$ws.Range('E28').Value = '  +0.08%  '

$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('E30').Value = '  +3.27%  '

$ws.Range('E31').Value = '  +0.26%  '

$ws.Range('D32').Value = '''7.97'
$ws.Range('E32').Value = '  +0.56%  '

$ws.Range('D33').Value = '''32.20'
$ws.Range('E33').Value = '  -2.47%  '

$ws.Range('D34').Value = '3.911.88'
$ws.Range('E34').Value = '  -0.96%  '

$ws.Range('D35').Value = '3.700.95'
$ws.Range('E35').Value = '  -1.17%  '

$ws.Range('E36').Value = '  -1.42%  '

$ws.Range('D37').Value = '''0.998'
$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('D38').Value = '''1.02'
$ws.Range('E38').Value = '  -0.27%  '

$ws.Range('D39').Value = '''5.87'
$ws.Range('E39').Value = '  +0.19%  '

$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('D41').Value = '''0.327'
$ws.Range('E41').Value = '  -1.04%  '

$ws.Range('D42').Value = '''438.28'
$ws.Range('E42').Value = '  -4.51%  '

$ws.Range('D43').Value = '''49.04'
$ws.Range('E43').Value = '  -0.25%  '

$ws.Range('D44').Value = '''1.99'
$ws.Range('E44').Value = '  -1.13%  '

$ws.Range('D45').Value = '''2.87'
$ws.Range('E45').Value = '  +0.18%  '

$ws.Range('D46').Value = '''8.54'
$ws.Range('E46').Value = '  +0.93%  '

$ws.Range('E47').Value = '  +0.01%  '

$ws.Range('D48').Value = '''40.67'
$ws.Range('E48').Value = '  +0.14%  '

$ws.Range('D49').Value = '2.829.54'
$ws.Range('E49').Value = '  -0.69%  '

$ws.Range('D50').Value = '''141.16'
$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('D51').Value = '''0.0357'
$ws.Range('E51').Value = '  +0.66%  '
